$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3.182878228561681,  1.65323645889881,   0.1529057820181812, 0.4998867070740569, 1, 5.488907176552729),
    @(3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(0.001754667048134761, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 1.551959061145647),
    @(1.505614041169197,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 1, 4.371470058157054),
    @(3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 1, 6.048734245549538),
    @(1.505614041169197,  1.65323645889881,   0.1529057820181812, 0.4998867070740569, 1, 3.811642989160245)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("F$row").Value = $vals[4]
    $ws.Range("G$row").Value = $vals[5]
}
